$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.075165666666667
$ws.Range("H2").Value = 9.225497
$ws.Range("I2").Value = 0.02641273658732285
$ws.Range("J2").Value = 0.02641273658732285
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 10.62858848795689
$ws.Range("R2").Value = 95.65729639161201
$ws.Range("S2").Value = 0.0002599418929028915
$ws.Range("T2").Value = 0.0002599418929028914
$ws.Range("G3").Value = 3.075165666666667
$ws.Range("H3").Value = 9.225497
$ws.Range("I3").Value = 0.02641273658732285
$ws.Range("J3").Value = 0.02641273658732285
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 927.4704222412959
$ws.Range("R3").Value = 8347.233800171662
$ws.Range("S3").Value = 0.02268301359508089
$ws.Range("T3").Value = 0.02268301359508089
$ws.Range("G4").Value = 3.075165666666667
$ws.Range("H4").Value = 9.225497
$ws.Range("I4").Value = 0.02641273658732285
$ws.Range("J4").Value = 0.02641273658732285
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 141.8735357980281
$ws.Range("R4").Value = 1276.861822182253
$ws.Range("S4").Value = 0.003469781099339063
$ws.Range("T4").Value = 0.003469781099339062
$ws.Range("I5").Value = 0.549422396165273
$ws.Range("J5").Value = 0.5494223961652731
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 221.089720695988
$ws.Range("R5").Value = 1989.807486263892
$ws.Range("S5").Value = 0.005407160185400509
$ws.Range("T5").Value = 0.005407160185400509
$ws.Range("I6").Value = 0.549422396165273
$ws.Range("J6").Value = 0.5494223961652731
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.4718388660885818
$ws.Range("T6").Value = 0.4718388660885818
$ws.Range("I7").Value = 0.549422396165273
$ws.Range("J7").Value = 0.5494223961652731
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.0721763698912908
$ws.Range("T7").Value = 0.0721763698912908
$ws.Range("I8").Value = 0.424164867247404
$ws.Range("J8").Value = 0.4241648672474041
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 170.6856012483511
$ws.Range("R8").Value = 1536.17041123516
$ws.Range("S8").Value = 0.004174433729374101
$ws.Range("T8").Value = 0.0041744337293741
$ws.Range("I9").Value = 0.424164867247404
$ws.Range("J9").Value = 0.4241648672474041
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.3642688601584147
$ws.Range("T9").Value = 0.3642688601584147
$ws.Range("I10").Value = 0.424164867247404
$ws.Range("J10").Value = 0.4241648672474041
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.05572157335961531
$ws.Range("T10").Value = 0.05572157335961531
